$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.388.42'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '1.687.07'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").Value = "'218.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("D6").Value = "'0.5538"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +9.11%  '
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").Value = "'0.2711"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("D9").Value = "'0.06493"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("D10").Value = "'22.14"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("D11").Value = "'0.07583"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("D12").Value = "'4.557"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.683.83'
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Value = "'0.5823"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = "'0.000008482"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = "'65.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '26.417.99'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = "'4.947"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = "'10.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").Value = "'191.27"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = "'6.248"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = "'149.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.20%  '
$ws.Range("D25").Value = "'0.1322"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +10.39%  '
$ws.Range("D26").Value = "'7.912"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.84%  '
$ws.Range("D27").Value = "'15.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("D28").Value = "'0.06340"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.21%  '
$ws.Range("D29").Value = "'1.396"
$ws.Range("D29").ClearFormats()
$ws.Range("D30").Value = "'1.329"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").Value = "'3.591"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.31%  '
$ws.Range("E32").Value = '  +2.11%  '
$ws.Range("D33").Value = "'1.677"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").Value = "'1.043"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.57%  '
$ws.Range("D35").Value = "'0.6246"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.94%  '
$ws.Range("D36").Value = "'2.409"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("D37").Value = "'2.716"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("D38").Value = "'6.242"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.01638"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.45%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.116.76'
$ws.Range("E40").Value = '  +2.05%  '
$ws.Range("D41").Value = "'0.8777"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("D42").Value = "'1.015"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = "'100.70"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '1.837.42'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.19%  '
$ws.Range("D46").Value = "'57.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").Value = "'8.227"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.91%  '
$ws.Range("D48").Value = "'1.006"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").Value = "'0.05286"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").Value = "'0.4300"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = "'6.091"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.95%  '
